$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.239.01"
$ws.Range("E2").Value = "  -6.06%  "
$ws.Range("D3").Value = "1.667.83"
$ws.Range("E3").Value = "  -4.15%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'217.50"
$ws.Range("E5").Value = "  -4.03%  "
$ws.Range("D6").Value = "'0.5068"
$ws.Range("E6").Value = "  -12.36%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'0.2648"
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").Value = "'0.06345"
$ws.Range("E9").Value = "  -4.30%  "
$ws.Range("D10").Value = "'21.45"
$ws.Range("E10").Value = "  -7.19%  "
$ws.Range("D11").Value = "'0.07366"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.669.22"
$ws.Range("E12").Value = "  -4.20%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.537"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "'0.5800"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("D15").Value = "1.894.59"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").Value = "'0.000008495"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "'64.56"
$ws.Range("E17").Value = "  -13.40%  "
$ws.Range("D18").Value = "26.293.04"
$ws.Range("E18").Value = "  -5.85%  "
$ws.Range("D19").Value = "'4.914"
$ws.Range("E19").Value = "  -7.49%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").Value = "'188.56"
$ws.Range("E22").Value = "  -8.16%  "
$ws.Range("D23").Value = "'6.174"
$ws.Range("E23").Value = "  -6.93%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "'143.22"
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").Value = "'7.669"
$ws.Range("E26").Value = "  -4.93%  "
$ws.Range("D27").Value = "'0.1172"
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("D28").Value = "'15.61"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").Value = "'0.05823"
$ws.Range("E29").Value = "  -6.03%  "
$ws.Range("D30").Value = "'1.271"
$ws.Range("E30").Value = "  -8.32%  "
$ws.Range("D31").Value = "'1.319"
$ws.Range("E31").Value = "  -5.34%  "
$ws.Range("D32").Value = "'3.531"
$ws.Range("E32").Value = "  -5.59%  "
$ws.Range("D33").Value = "'3.505"
$ws.Range("D34").Value = "'1.632"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("E36").Value = "  -6.73%  "
$ws.Range("D37").Value = "'2.360"
$ws.Range("D38").Value = "'2.644"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("D40").Value = "'6.007"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "1.071.96"
$ws.Range("E41").Value = "  -4.56%  "
$ws.Range("D42").Value = "'0.8656"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'99.52"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "1.816.53"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").Value = "'0.00000000110"
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("D47").Value = "'55.67"
$ws.Range("E47").Value = "  -6.25%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "'8.073"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "'0.4293"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").Value = "'0.05178"
$ws.Range("E51").Value = "  -3.77%  "
